# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund holdings detail) right before
#    the "总计" (total) summary sheet, duplicating an existing quarter
#    sheet so that all sheet-level formatting/styles (borders, bold
#    header, page margins, etc.) are preserved exactly.
# 2) Populate the new sheet with the 2022-Q1 fund holding data.
# 3) Insert a new row at the top of the "总计" data (row 2) for 2022-Q1
#    and renumber the existing rows' index column.

$wb = $excel.ActiveWorkbook

# Helper: write $value into $range as literal text, even when $value
# looks like a number (e.g. "010966", "59.70"), and do so WITHOUT leaving
# any NumberFormat/style residue on $range (matching the source file's
# convention of plain, unstyled inlineStr cells for this kind of data).
# This is done by staging the text in a scratch cell (forced to text via
# NumberFormat "@"), copying only its value+type into the destination via
# PasteSpecial(xlPasteValues), then wiping the scratch cell.
function Set-TextValue($range, $value) {
    $scratch = $range.Worksheet.Range("ZZ9000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# Duplicate the template sheet, placing the copy immediately before "总计".
# Using Worksheet.Copy (rather than Worksheets.Add + PasteSpecial) preserves
# the sheetPr/pageMargins/sheetFormatPr and exact cell styles of the source.
$templateSheet.Copy($totalSheet)

# Re-fetch the "总计" sheet: the previously held COM reference's cached
# Index does not update after the copy operation inserts a sheet before it.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# --- Header row (cells already carry the right bold/border style from the
#     duplicated template, so a plain text assignment is enough) ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Data rows ---
# Column A: plain numeric row index (0-based), keeps the template's "s=2"
#           style, so a direct numeric assignment is fine.
# Columns B, D, E, F, G: numeric-looking text (fund code / scale / weight /
#           market value), must be forced to text so leading zeros and the
#           original "t=inlineStr" representation are preserved.
# Column C: fund name, plain non-numeric text - direct assignment is fine.
# Column H: genuine integer rank - direct numeric assignment is fine.
$rows = @(
    @{ idx = 0; code = "010966"; name = "富国成长领航混合";              scale = "59.70"; pos = "88.71"; pct = "7.08"; mv = "4.2268"; rank = 1 },
    @{ idx = 1; code = "010662"; name = "富国均衡优选混合";              scale = "49.13"; pos = "91.70"; pct = "7.72"; mv = "3.7928"; rank = 1 },
    @{ idx = 2; code = "001985"; name = "富国低碳新经济混合A";           scale = "39.35"; pos = "93.88"; pct = "7.36"; mv = "2.8962"; rank = 2 },
    @{ idx = 3; code = "005368"; name = "富国清洁能源产业灵活配置混合A"; scale = "23.52"; pos = "88.60"; pct = "2.89"; mv = "0.6797"; rank = 9 },
    @{ idx = 4; code = "011127"; name = "富国清洁能源产业灵活配置混合C"; scale = "1.61";  pos = "88.60"; pct = "2.89"; mv = "0.0465"; rank = 9 },
    @{ idx = 5; code = "011306"; name = "富国低碳新经济混合C";           scale = "0.30";  pos = "93.88"; pct = "7.36"; mv = "0.0221"; rank = 2 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row.idx
    Set-TextValue $newSheet.Range("B$r") $row.code
    $newSheet.Range("C$r").Value = $row.name
    Set-TextValue $newSheet.Range("D$r") $row.scale
    Set-TextValue $newSheet.Range("E$r") $row.pos
    Set-TextValue $newSheet.Range("F$r") $row.pct
    Set-TextValue $newSheet.Range("G$r") $row.mv
    $newSheet.Range("H$r").Value = $row.rank
    $r = $r + 1
}

# --- Update the "总计" (total) sheet: insert the 2022-Q1 summary row ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# The freshly inserted row's A2 cell does not inherit the bold/bordered
# index-column style used by the rest of column A, so copy it explicitly
# from the row below (format only - the value we set next is untouched).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 11.66

# Renumber the index column (A) for the rows that shifted down by one.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
